$wb = $excel.ActiveWorkbook

# This script applies refreshed market-price-derived values (columns H-N:
# currentAveragePrice, currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, as produced by the scheduled
# market-data refresh run.

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 195  # H9
$ws.Cells.Item(9, 9).Value = 102.55  # I9
$ws.Cells.Item(9, 10).Value = 657.25  # J9
$ws.Cells.Item(9, 11).Value = 102.55  # K9
$ws.Cells.Item(9, 12).Value = 657.25  # L9
$ws.Cells.Item(9, 13).Value = 66.45  # M9
$ws.Cells.Item(9, 14).Value = -995.25  # N9
$ws.Cells.Item(19, 8).Value = 2371.1052  # H19
$ws.Cells.Item(19, 9).Value = 1523.2222  # I19
$ws.Cells.Item(19, 10).Value = 3134.2  # J19
$ws.Cells.Item(19, 11).Value = 1523.2222  # K19
$ws.Cells.Item(19, 12).Value = 3134.2  # L19
$ws.Cells.Item(19, 13).Value = -1348.2222  # M19
$ws.Cells.Item(19, 14).Value = -3484.2  # N19
$ws.Cells.Item(31, 8).Value = 199.57143  # H31
$ws.Cells.Item(31, 9).Value = 199.57143  # I31
$ws.Cells.Item(31, 10).Value = 0  # J31
$ws.Cells.Item(31, 11).Value = 598.71429  # K31
$ws.Cells.Item(31, 12).Value = 0  # L31
$ws.Cells.Item(31, 13).Value = -368.71429  # M31
$ws.Cells.Item(31, 14).ClearContents()  # N31
$ws.Cells.Item(33, 8).Value = 302.6875  # H33
$ws.Cells.Item(33, 9).Value = 256.2  # I33
$ws.Cells.Item(33, 11).Value = 256.2  # K33
$ws.Cells.Item(33, 13).Value = -27.19999999999999  # M33
$ws.Cells.Item(113, 8).Value = 4934.857  # H113
$ws.Cells.Item(113, 9).Value = 4867.6665  # I113
$ws.Cells.Item(113, 11).Value = 4867.6665  # K113
$ws.Cells.Item(113, 13).Value = -1613.6665  # M113
$ws.Cells.Item(125, 8).Value = 11070.777  # H125
$ws.Cells.Item(125, 10).Value = 8842.5  # J125
$ws.Cells.Item(125, 12).Value = 79582.5  # L125
$ws.Cells.Item(125, 14).Value = -84502.5  # N125
$ws.Cells.Item(135, 8).Value = 1229.1177  # H135
$ws.Cells.Item(135, 9).Value = 1193  # I135
$ws.Cells.Item(135, 11).Value = 10737  # K135
$ws.Cells.Item(135, 13).Value = -8202  # M135

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2684.7646  # H2
$ws.Cells.Item(2, 9).Value = 2216.7856  # I2
$ws.Cells.Item(2, 11).Value = 2216.7856  # K2
$ws.Cells.Item(2, 13).Value = -2103.7856  # M2
$ws.Cells.Item(5, 8).Value = 349.6  # H5
$ws.Cells.Item(5, 9).Value = 312  # I5
$ws.Cells.Item(5, 11).Value = 312  # K5
$ws.Cells.Item(5, 13).Value = -200  # M5
$ws.Cells.Item(45, 8).Value = 2617.5  # H45
$ws.Cells.Item(45, 9).Value = 2727.5833  # I45
$ws.Cells.Item(45, 10).Value = 1957  # J45
$ws.Cells.Item(45, 11).Value = 2727.5833  # K45
$ws.Cells.Item(45, 12).Value = 1957  # L45
$ws.Cells.Item(45, 13).Value = -2350.5833  # M45
$ws.Cells.Item(45, 14).Value = -2711  # N45
$ws.Cells.Item(58, 8).Value = 0  # H58
$ws.Cells.Item(58, 10).Value = 0  # J58
$ws.Cells.Item(58, 12).Value = 0  # L58
$ws.Cells.Item(58, 14).ClearContents()  # N58
$ws.Cells.Item(63, 8).Value = 2775.4443  # H63
$ws.Cells.Item(63, 9).Value = 2775.4443  # I63
$ws.Cells.Item(63, 11).Value = 2775.4443  # K63
$ws.Cells.Item(63, 13).Value = -2089.4443  # M63
$ws.Cells.Item(66, 8).Value = 2775.4443  # H66
$ws.Cells.Item(66, 9).Value = 2775.4443  # I66
$ws.Cells.Item(66, 11).Value = 13877.2215  # K66
$ws.Cells.Item(66, 13).Value = -10445.2215  # M66
$ws.Cells.Item(88, 8).Value = 1724.0968  # H88
$ws.Cells.Item(88, 9).Value = 1850.5  # I88
$ws.Cells.Item(88, 10).Value = 1549.0769  # J88
$ws.Cells.Item(88, 11).Value = 1850.5  # K88
$ws.Cells.Item(88, 12).Value = 1549.0769  # L88
$ws.Cells.Item(88, 13).Value = -1444.5  # M88
$ws.Cells.Item(88, 14).Value = -2361.0769  # N88
$ws.Cells.Item(91, 8).Value = 1724.0968  # H91
$ws.Cells.Item(91, 9).Value = 1850.5  # I91
$ws.Cells.Item(91, 10).Value = 1549.0769  # J91
$ws.Cells.Item(91, 11).Value = 1850.5  # K91
$ws.Cells.Item(91, 12).Value = 1549.0769  # L91
$ws.Cells.Item(91, 13).Value = -446.5  # M91
$ws.Cells.Item(91, 14).Value = -4357.0769  # N91
$ws.Cells.Item(116, 8).Value = 2684.7646  # H116
$ws.Cells.Item(116, 9).Value = 2216.7856  # I116
$ws.Cells.Item(116, 11).Value = 2216.7856  # K116
$ws.Cells.Item(116, 13).Value = 77.21439999999984  # M116
$ws.Cells.Item(138, 8).Value = 64429  # H138
$ws.Cells.Item(138, 10).Value = 64429  # J138
$ws.Cells.Item(138, 12).Value = 64429  # L138
$ws.Cells.Item(138, 14).Value = -74709  # N138

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2684.7646  # H3
$ws.Cells.Item(3, 9).Value = 2216.7856  # I3
$ws.Cells.Item(3, 11).Value = 2216.7856  # K3
$ws.Cells.Item(3, 13).Value = -2102.7856  # M3
$ws.Cells.Item(4, 8).Value = 349.6  # H4
$ws.Cells.Item(4, 9).Value = 312  # I4
$ws.Cells.Item(4, 11).Value = 312  # K4
$ws.Cells.Item(4, 13).Value = -197  # M4
$ws.Cells.Item(11, 8).Value = 1378.0714  # H11
$ws.Cells.Item(11, 9).Value = 1160.4  # I11
$ws.Cells.Item(11, 10).Value = 1499  # J11
$ws.Cells.Item(11, 11).Value = 1160.4  # K11
$ws.Cells.Item(11, 12).Value = 1499  # L11
$ws.Cells.Item(11, 13).Value = -1020.4  # M11
$ws.Cells.Item(11, 14).Value = -1779  # N11
$ws.Cells.Item(37, 8).Value = 1200  # H37
$ws.Cells.Item(37, 9).Value = 1200  # I37
$ws.Cells.Item(37, 11).Value = 1200  # K37
$ws.Cells.Item(37, 13).Value = -1063  # M37
$ws.Cells.Item(92, 8).Value = 22000  # H92
$ws.Cells.Item(92, 10).Value = 22000  # J92
$ws.Cells.Item(92, 12).Value = 22000  # L92
$ws.Cells.Item(92, 14).Value = -26992  # N92

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 512.6667  # H7
$ws.Cells.Item(7, 9).Value = 564.25  # I7
$ws.Cells.Item(7, 10).Value = 100  # J7
$ws.Cells.Item(7, 11).Value = 564.25  # K7
$ws.Cells.Item(7, 12).Value = 100  # L7
$ws.Cells.Item(7, 13).Value = -451.25  # M7
$ws.Cells.Item(7, 14).Value = -326  # N7
$ws.Cells.Item(22, 8).Value = 23  # H22
$ws.Cells.Item(22, 9).Value = 23  # I22
$ws.Cells.Item(22, 11).Value = 23  # K22
$ws.Cells.Item(22, 13).Value = 327  # M22
$ws.Cells.Item(31, 8).Value = 592032.5  # H31
$ws.Cells.Item(31, 9).Value = 1253626.1  # I31
$ws.Cells.Item(31, 10).Value = 3949.2222  # J31
$ws.Cells.Item(31, 11).Value = 1253626.1  # K31
$ws.Cells.Item(31, 12).Value = 3949.2222  # L31
$ws.Cells.Item(31, 13).Value = -1253331.1  # M31
$ws.Cells.Item(31, 14).Value = -4539.2222  # N31
$ws.Cells.Item(34, 8).Value = 592032.5  # H34
$ws.Cells.Item(34, 9).Value = 1253626.1  # I34
$ws.Cells.Item(34, 10).Value = 3949.2222  # J34
$ws.Cells.Item(34, 11).Value = 1253626.1  # K34
$ws.Cells.Item(34, 12).Value = 3949.2222  # L34
$ws.Cells.Item(34, 13).Value = -1253424.1  # M34
$ws.Cells.Item(34, 14).Value = -4353.2222  # N34
$ws.Cells.Item(99, 8).Value = 786381.1  # H99
$ws.Cells.Item(99, 10).Value = 19299.3  # J99
$ws.Cells.Item(99, 12).Value = 19299.3  # L99
$ws.Cells.Item(99, 14).Value = -22295.3  # N99
$ws.Cells.Item(126, 8).Value = 786381.1  # H126
$ws.Cells.Item(126, 10).Value = 19299.3  # J126
$ws.Cells.Item(126, 12).Value = 57897.89999999999  # L126
$ws.Cells.Item(126, 14).Value = -62837.89999999999  # N126

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 463932.25  # H4
$ws.Cells.Item(4, 9).Value = 345128.28  # I4
$ws.Cells.Item(4, 10).Value = 630257.8  # J4
$ws.Cells.Item(4, 11).Value = 1035384.84  # K4
$ws.Cells.Item(4, 12).Value = 1890773.4  # L4
$ws.Cells.Item(4, 13).Value = -1035272.84  # M4
$ws.Cells.Item(4, 14).Value = -1890997.4  # N4
$ws.Cells.Item(55, 8).Value = 8752.5  # H55
$ws.Cells.Item(55, 9).Value = 7000  # I55
$ws.Cells.Item(55, 10).Value = 8947.223  # J55
$ws.Cells.Item(55, 11).Value = 21000  # K55
$ws.Cells.Item(55, 12).Value = 26841.669  # L55
$ws.Cells.Item(55, 13).Value = -20823  # M55
$ws.Cells.Item(55, 14).Value = -27195.669  # N55
$ws.Cells.Item(56, 8).Value = 5630  # H56
$ws.Cells.Item(56, 9).Value = 5630  # I56
$ws.Cells.Item(56, 11).Value = 5630  # K56
$ws.Cells.Item(56, 13).Value = -5100  # M56
$ws.Cells.Item(81, 8).Value = 399.5  # H81
$ws.Cells.Item(81, 9).Value = 399.5  # I81
$ws.Cells.Item(81, 11).Value = 1198.5  # K81
$ws.Cells.Item(81, 13).Value = -75.5  # M81
$ws.Cells.Item(84, 8).Value = 399.5  # H84
$ws.Cells.Item(84, 9).Value = 399.5  # I84
$ws.Cells.Item(84, 11).Value = 3595.5  # K84
$ws.Cells.Item(84, 13).Value = 2020.5  # M84
$ws.Cells.Item(88, 8).Value = 7007  # H88
$ws.Cells.Item(88, 10).Value = 7007  # J88
$ws.Cells.Item(88, 12).Value = 21021  # L88
$ws.Cells.Item(88, 14).Value = -21877  # N88
$ws.Cells.Item(91, 8).Value = 7007  # H91
$ws.Cells.Item(91, 10).Value = 7007  # J91
$ws.Cells.Item(91, 12).Value = 21021  # L91
$ws.Cells.Item(91, 14).Value = -23985  # N91
$ws.Cells.Item(113, 8).Value = 835.1667  # H113
$ws.Cells.Item(113, 10).Value = 820.1818  # J113
$ws.Cells.Item(113, 12).Value = 2460.5454  # L113
$ws.Cells.Item(113, 14).Value = -6800.5454  # N113
$ws.Cells.Item(131, 8).Value = 1184.4375  # H131
$ws.Cells.Item(131, 9).Value = 621.9  # I131
$ws.Cells.Item(131, 10).Value = 1440.1364  # J131
$ws.Cells.Item(131, 11).Value = 1865.7  # K131
$ws.Cells.Item(131, 12).Value = 4320.4092  # L131
$ws.Cells.Item(131, 13).Value = 3174.3  # M131
$ws.Cells.Item(131, 14).Value = -14400.4092  # N131

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 205.4  # H2
$ws.Cells.Item(2, 9).Value = 172.66667  # I2
$ws.Cells.Item(2, 10).Value = 500  # J2
$ws.Cells.Item(2, 11).Value = 172.66667  # K2
$ws.Cells.Item(2, 12).Value = 500  # L2
$ws.Cells.Item(2, 13).Value = -59.66667000000001  # M2
$ws.Cells.Item(2, 14).Value = -726  # N2
$ws.Cells.Item(113, 8).Value = 3195.875  # H113
$ws.Cells.Item(113, 9).Value = 3206.3333  # I113
$ws.Cells.Item(113, 10).Value = 3189.6  # J113
$ws.Cells.Item(113, 11).Value = 3206.3333  # K113
$ws.Cells.Item(113, 12).Value = 3189.6  # L113
$ws.Cells.Item(113, 13).Value = -1036.3333  # M113
$ws.Cells.Item(113, 14).Value = -7529.6  # N113
$ws.Cells.Item(132, 8).Value = 20581.379  # H132
$ws.Cells.Item(132, 9).Value = 23848.857  # I132
$ws.Cells.Item(132, 10).Value = 2791.7778  # J132
$ws.Cells.Item(132, 11).Value = 71546.571  # K132
$ws.Cells.Item(132, 12).Value = 8375.3334  # L132
$ws.Cells.Item(132, 13).Value = -69016.571  # M132
$ws.Cells.Item(132, 14).Value = -13435.3334  # N132
$ws.Cells.Item(141, 8).Value = 42160.5  # H141
$ws.Cells.Item(141, 10).Value = 44000  # J141
$ws.Cells.Item(141, 12).Value = 44000  # L141
$ws.Cells.Item(141, 14).Value = -54360  # N141

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 4333.6665  # H122
$ws.Cells.Item(122, 9).Value = 4250.25  # I122
$ws.Cells.Item(122, 11).Value = 12750.75  # K122
$ws.Cells.Item(122, 13).Value = -10300.75  # M122
$ws.Cells.Item(136, 8).Value = 3235.0908  # H136
$ws.Cells.Item(136, 9).Value = 1699.125  # I136
$ws.Cells.Item(136, 11).Value = 5097.375  # K136
$ws.Cells.Item(136, 13).Value = -2547.375  # M136

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(23, 8).Value = 1755  # H23
$ws.Cells.Item(23, 9).Value = 1206  # I23
$ws.Cells.Item(23, 11).Value = 1206  # K23
$ws.Cells.Item(23, 13).Value = -977  # M23
$ws.Cells.Item(126, 8).Value = 3393.0527  # H126
$ws.Cells.Item(126, 9).Value = 3091.9375  # I126
$ws.Cells.Item(126, 11).Value = 9275.8125  # K126
$ws.Cells.Item(126, 13).Value = -6805.8125  # M126
$ws.Cells.Item(140, 8).Value = 55071.75  # H140
$ws.Cells.Item(140, 10).Value = 55071.75  # J140
$ws.Cells.Item(140, 12).Value = 55071.75  # L140
$ws.Cells.Item(140, 14).Value = -65431.75  # N140
$ws.Cells.Item(141, 8).Value = 68301.55499999999  # H141
$ws.Cells.Item(141, 10).Value = 68301.55499999999  # J141
$ws.Cells.Item(141, 12).Value = 68301.55499999999  # L141
$ws.Cells.Item(141, 14).Value = -78661.55499999999  # N141
